# Atualização de bases das ligas, do dia: 20-06-2024 às 20:11
# Swap the match-data columns (B:AD) between specific row pairs while
# leaving column A (the row index) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(49, 50),
    @(177, 178),
    @(190, 191),
    @(192, 194),
    @(198, 199)
)

foreach ($pair in $pairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    $rangeA = $ws.Range("B$rowA`:AD$rowA")
    $rangeB = $ws.Range("B$rowB`:AD$rowB")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}
